$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.450.41"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "2.943.74"
$ws.Range("E3").Value = "  -3.91%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'493.20"
$ws.Range("E5").Value = "  -6.70%  "
$ws.Range("D6").Value = "'133.48"
$ws.Range("E6").Value = "  -7.13%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.422"
$ws.Range("E8").Value = "  -6.01%  "
$ws.Range("D9").Value = "'7.09"
$ws.Range("E9").Value = "  -7.46%  "
$ws.Range("E10").Value = "  -7.63%  "
$ws.Range("D11").Value = "'0.349"
$ws.Range("E11").Value = "  -6.15%  "
$ws.Range("D12").Value = "3.456.47"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("E13").Value = "  -3.69%  "
$ws.Range("D14").Value = "'25.77"
$ws.Range("E14").Value = "  -6.09%  "
$ws.Range("D15").Value = "'0.0000155"
$ws.Range("E15").Value = "  -10.15%  "
$ws.Range("D16").Value = "56.529.88"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "2.946.76"
$ws.Range("E17").Value = "  -3.88%  "
$ws.Range("D18").Value = "'5.93"
$ws.Range("E18").Value = "  -4.72%  "
$ws.Range("D19").Value = "'12.36"
$ws.Range("E19").Value = "  -6.47%  "
$ws.Range("D20").Value = "'7.70"
$ws.Range("E20").Value = "  -6.38%  "
$ws.Range("D21").Value = "'314.41"
$ws.Range("E21").Value = "  -8.16%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'5.70"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "'0.480"
$ws.Range("E24").Value = "  -5.12%  "
$ws.Range("D25").Value = "'62.39"
$ws.Range("E25").Value = "  -4.72%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = "  -5.78%  "
$ws.Range("D28").Value = "0.0₃0849"
$ws.Range("E28").Value = "  -13.38%  "
$ws.Range("D29").Value = "'6.38"
$ws.Range("E29").Value = "  -9.13%  "
$ws.Range("D30").Value = "'6.95"
$ws.Range("E30").Value = "  -7.16%  "
$ws.Range("D31").Value = "'1.73"
$ws.Range("E31").Value = "  -6.90%  "
$ws.Range("D32").Value = "'19.81"
$ws.Range("E32").Value = "  -6.75%  "
$ws.Range("E33").Value = "  -9.78%  "
$ws.Range("D34").Value = "'151.42"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").Value = "'4.42"
$ws.Range("E35").Value = "  -8.53%  "
$ws.Range("D36").Value = "'5.63"
$ws.Range("E36").Value = "  -5.98%  "
$ws.Range("D37").Value = "'1.19"
$ws.Range("E37").Value = "  -10.61%  "
$ws.Range("D38").Value = "'23.50"
$ws.Range("E38").Value = "  -10.02%  "
$ws.Range("D39").Value = "'0.0647"
$ws.Range("E39").Value = "  -7.63%  "
$ws.Range("D40").Value = "2.976.17"
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("D41").Value = "'37.21"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'0.635"
$ws.Range("E43").Value = "  -4.86%  "
$ws.Range("D44").Value = "'3.65"
$ws.Range("E44").Value = "  -8.77%  "
$ws.Range("D45").Value = "2.131.65"
$ws.Range("E45").Value = "  -8.95%  "
$ws.Range("E46").Value = "  -10.11%  "
$ws.Range("D47").Value = "'5.81"
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("D48").Value = "'0.912"
$ws.Range("E48").Value = "  -12.70%  "
$ws.Range("D49").Value = "'0.0228"
$ws.Range("E49").Value = "  -7.06%  "
$ws.Range("D50").Value = "'18.72"
$ws.Range("E50").Value = "  -8.18%  "
$ws.Range("D51").Value = "'1.71"
$ws.Range("E51").Value = "  -15.44%  "
